$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '24.744.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.701.95'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.22'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4057'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.516'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.002'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.10'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08888'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.658'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.70'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.153'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001327'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.693.12'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.50'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07099'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.133'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.01%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '24.749.71'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.142'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  +1.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.80'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.269'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +22.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '164.81'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.45%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '138.88'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.127'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.959'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09085'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.074'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02986'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2791'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.39%  '
$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.32%  '
$ws.Range('B38').Value = 'WEMIXTOKEN'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.959'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09264'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7812'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.470'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.17'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.46%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.653'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.39%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7255'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.213'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.358'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.51%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07992'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('E51').Value = '  +2.40%  '
